$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: round the Ost/Nord coordinate values
$ws.Range("Q12").Value = 565669
$ws.Range("R12").Value = 6699889

# Row 12: remove the Starttid (Z12) and Sluttid (AB12) values entirely
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()

# Row 13: round the Ost/Nord coordinate values
$ws.Range("Q13").Value = 565683
$ws.Range("R13").Value = 6699911

# Row 13: remove the Starttid (Z13) and Sluttid (AB13) values entirely
$ws.Range("Z13").ClearContents()
$ws.Range("AB13").ClearContents()
